# Update cryptocurrency price/volume snapshot (and restore the Quant/NEARProtocol
# row ordering) to match the refreshed GitHub Actions data pull.
#
# Column D ("Price") values are plain decimal-looking text (e.g. "1.000",
# "0.00001080") that must stay TEXT cells -- exactly like the source data --
# rather than being auto-coerced into numbers (which would strip trailing
# zeros / switch to scientific notation). We force that with a leading
# apostrophe, same as typing '1.000 directly into Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; B = $null; C = $null; D = "28.276.72"; E = "  +2.89%  " },
    @{ Row = 3; B = $null; C = $null; D = "1.816.59"; E = "  +4.10%  " },
    @{ Row = 4; B = $null; C = $null; D = "1.003"; E = "  +0.06%  " },
    @{ Row = 5; B = $null; C = $null; D = "328.29"; E = "  +2.06%  " },
    @{ Row = 6; B = $null; C = $null; D = "1.001"; E = "  +0.04%  " },
    @{ Row = 7; B = $null; C = $null; D = "0.4347"; E = "  +3.34%  " },
    @{ Row = 8; B = $null; C = $null; D = "0.3671"; E = "  +2.61%  " },
    @{ Row = 9; B = $null; C = $null; D = "44.93"; E = "  -1.14%  " },
    @{ Row = 10; B = $null; C = $null; D = "0.07679"; E = "  +3.73%  " },
    @{ Row = 11; B = $null; C = $null; D = $null; E = "  +2.77%  " },
    @{ Row = 12; B = $null; C = $null; D = $null; E = "  +0.05%  " },
    @{ Row = 13; B = $null; C = $null; D = "22.15"; E = "  +3.59%  " },
    @{ Row = 14; B = $null; C = $null; D = "6.300"; E = "  +3.34%  " },
    @{ Row = 15; B = $null; C = $null; D = "7.536"; E = "  +5.06%  " },
    @{ Row = 16; B = $null; C = $null; D = "1.830.98"; E = "  +4.95%  " },
    @{ Row = 17; B = $null; C = $null; D = "93.39"; E = "  +6.40%  " },
    @{ Row = 18; B = $null; C = $null; D = "0.00001080"; E = "  +1.69%  " },
    @{ Row = 19; B = $null; C = $null; D = $null; E = "  +7.13%  " },
    @{ Row = 20; B = $null; C = $null; D = $null; E = "  -0.02%  " },
    @{ Row = 21; B = $null; C = $null; D = "17.52"; E = "  +4.15%  " },
    @{ Row = 22; B = $null; C = $null; D = "6.264"; E = "  +3.13%  " },
    @{ Row = 23; B = $null; C = $null; D = "28.309.34"; E = "  +2.97%  " },
    @{ Row = 24; B = $null; C = $null; D = "11.63"; E = "  +1.66%  " },
    @{ Row = 25; B = $null; C = $null; D = "2.026"; E = "  -13.20%  " },
    @{ Row = 26; B = $null; C = $null; D = "162.41"; E = "  +6.57%  " },
    @{ Row = 27; B = $null; C = $null; D = "20.73"; E = "  +2.10%  " },
    @{ Row = 28; B = $null; C = $null; D = "2.035.94"; E = "  +4.84%  " },
    @{ Row = 29; B = $null; C = $null; D = "2.300"; E = "  -2.78%  " },
    @{ Row = 30; B = $null; C = $null; D = "128.82"; E = "  +2.64%  " },
    @{ Row = 31; B = $null; C = $null; D = $null; E = "  +1.26%  " },
    @{ Row = 32; B = $null; C = $null; D = "5.947"; E = "  +5.25%  " },
    @{ Row = 33; B = $null; C = $null; D = "0.09192"; E = "  +0.88%  " },
    @{ Row = 35; B = $null; C = $null; D = $null; E = "  +3.01%  " },
    @{ Row = 36; B = $null; C = $null; D = "0.02344"; E = "  +2.37%  " },
    @{ Row = 37; B = $null; C = $null; D = "0.2179"; E = "  +2.30%  " },
    @{ Row = 38; B = $null; C = $null; D = "5.201"; E = "  +2.73%  " },
    @{ Row = 39; B = $null; C = $null; D = "0.6571"; E = "  +3.39%  " },
    @{ Row = 40; B = $null; C = $null; D = "0.06201"; E = "  +2.66%  " },
    @{ Row = 41; B = $null; C = $null; D = "1.192"; E = "  +0.44%  " },
    @{ Row = 42; B = $null; C = $null; D = "8.126"; E = "  +3.13%  " },
    @{ Row = 43; B = $null; C = $null; D = "1.428"; E = "  -0.68%  " },
    @{ Row = 44; B = $null; C = $null; D = "1.000"; E = "  +0.04%  " },
    @{ Row = 45; B = $null; C = $null; D = "13.83"; E = "  +1.52%  " },
    @{ Row = 46; B = $null; C = $null; D = "0.6120"; E = "  +4.86%  " },
    @{ Row = 47; B = $null; C = $null; D = "3.757"; E = "  +1.39%  " },
    @{ Row = 48; B = "Quant"; C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D = "125.70"; E = "  +0.78%  " },
    @{ Row = 49; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "2.020"; E = "  +4.15%  " },
    @{ Row = 50; B = $null; C = $null; D = "1.158"; E = "  +4.17%  " },
    @{ Row = 51; B = $null; C = $null; D = "0.07004"; E = "  +2.40%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.B) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($null -ne $u.C) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($null -ne $u.D) {
        # Leading "'" forces text storage so numeric-looking strings keep
        # their exact formatting (trailing zeros, leading zeros, etc.)
        $ws.Range("D$row").Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}
